$wb = $excel.ActiveWorkbook

# Fill in B1 (TxHash values) and select P49 as the last active cell on that sheet
$ws1 = $wb.Worksheets.Item("B1")
$ws1.Activate()
$ws1.Range("A2").Value = "FB17469645ACDBDA2CD7C7EF27063C1DFA88C4CE0CFBECB231D49F4E4FBB6A33"
$ws1.Range("A3").Value = "92945CAA5FAE911EA7270912CDA3CBC3412FB4481318324EB7162B8D7BB86DF9"
$ws1.Range("P49").Select()

# Fill in B2 (TxHash values) and leave it as the final active sheet with A3 selected
$ws2 = $wb.Worksheets.Item("B2")
$ws2.Activate()
$ws2.Range("A2").Value = "F901B7318DD925DA53168062BB795E6A306BCB2B26327BCF55987998A67DF384"
$ws2.Range("A3").Value = "0D7643BA695D83810B41A4AC279119FE67822A8DB6203B2C6DB7998B3E0DA1A6"
$ws2.Range("A3").Select()
